# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates to ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets
# as captured by the commit diff (scheduled runner data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(46, 8).Value = 6600  # H46: was 86730
$ws.Cells.Item(46, 10).Value = 5000  # J46: was 102116
$ws.Cells.Item(46, 12).Value = 15000  # L46: was 306348
$ws.Cells.Item(46, 14).Value = -15238  # N46: was -306586
$ws.Cells.Item(60, 8).Value = 6600  # H60: was 86730
$ws.Cells.Item(60, 10).Value = 5000  # J60: was 102116
$ws.Cells.Item(60, 12).Value = 15000  # L60: was 306348
$ws.Cells.Item(60, 14).Value = -15968  # N60: was -307316
$ws.Cells.Item(62, 8).Value = 7078.3335  # H62: was 6501.1
$ws.Cells.Item(62, 10).Value = 14400  # J62: was 11126.5
$ws.Cells.Item(62, 12).Value = 14400  # L62: was 11126.5
$ws.Cells.Item(62, 14).Value = -15648  # N62: was -12374.5
$ws.Cells.Item(65, 8).Value = 7078.3335  # H65: was 6501.1
$ws.Cells.Item(65, 10).Value = 14400  # J65: was 11126.5
$ws.Cells.Item(65, 12).Value = 72000  # L65: was 55632.5
$ws.Cells.Item(65, 14).Value = -78240  # N65: was -61872.5
$ws.Cells.Item(86, 8).Value = 3980  # H86: was 169834.5
$ws.Cells.Item(86, 9).Value = 0  # I86: was 1000003
$ws.Cells.Item(86, 10).Value = 3980  # J86: was 3800.8
$ws.Cells.Item(86, 11).Value = 0  # K86: was 1000003
$ws.Cells.Item(86, 12).Value = 3980  # L86: was 3800.8
$ws.Cells.Item(86, 13).ClearContents()  # M86: was -998880
$ws.Cells.Item(86, 14).Value = -6226  # N86: was -6046.8
$ws.Cells.Item(89, 8).Value = 3980  # H89: was 169834.5
$ws.Cells.Item(89, 9).Value = 0  # I89: was 1000003
$ws.Cells.Item(89, 10).Value = 3980  # J89: was 3800.8
$ws.Cells.Item(89, 11).Value = 0  # K89: was 5000015
$ws.Cells.Item(89, 12).Value = 19900  # L89: was 19004
$ws.Cells.Item(89, 13).ClearContents()  # M89: was -4994399
$ws.Cells.Item(89, 14).Value = -31132  # N89: was -30236
$ws.Cells.Item(98, 8).Value = 1023.63635  # H98: was 901.7
$ws.Cells.Item(98, 9).Value = 998.9286  # I98: was 732
$ws.Cells.Item(98, 10).Value = 1066.875  # J98: was 1071.4
$ws.Cells.Item(98, 11).Value = 998.9286  # K98: was 732
$ws.Cells.Item(98, 12).Value = 1066.875  # L98: was 1071.4
$ws.Cells.Item(98, 13).Value = 499.0714  # M98: was 766
$ws.Cells.Item(98, 14).Value = -4062.875  # N98: was -4067.4
$ws.Cells.Item(106, 8).Value = 3480.4  # H106: was 3538.125
$ws.Cells.Item(106, 9).Value = 3828.2856  # I106: was 4059.8
$ws.Cells.Item(106, 11).Value = 3828.2856  # K106: was 4059.8
$ws.Cells.Item(106, 13).Value = -3197.2856  # M106: was -3428.8
$ws.Cells.Item(122, 8).Value = 1023.63635  # H122: was 901.7
$ws.Cells.Item(122, 9).Value = 998.9286  # I122: was 732
$ws.Cells.Item(122, 10).Value = 1066.875  # J122: was 1071.4
$ws.Cells.Item(122, 11).Value = 2996.7858  # K122: was 2196
$ws.Cells.Item(122, 12).Value = 3200.625  # L122: was 3214.2
$ws.Cells.Item(122, 13).Value = -546.7857999999997  # M122: was 254
$ws.Cells.Item(122, 14).Value = -8100.625  # N122: was -8114.200000000001
$ws.Cells.Item(129, 8).Value = 1013.7692  # H129: was 1006.01666
$ws.Cells.Item(129, 9).Value = 377.75  # I129: was 398.2
$ws.Cells.Item(129, 10).Value = 1066.7709  # J129: was 1061.2727
$ws.Cells.Item(129, 11).Value = 1133.25  # K129: was 1194.6
$ws.Cells.Item(129, 12).Value = 3200.3127  # L129: was 3183.8181
$ws.Cells.Item(129, 13).Value = 3866.75  # M129: was 3805.4
$ws.Cells.Item(129, 14).Value = -13200.3127  # N129: was -13183.8181
$ws.Cells.Item(135, 8).Value = 100001520  # H135: was 88236620
$ws.Cells.Item(135, 9).Value = 45456556  # I135: was 31251414
$ws.Cells.Item(135, 10).Value = 250000180  # J135: was 1000000000
$ws.Cells.Item(135, 11).Value = 409109004  # K135: was 281262726
$ws.Cells.Item(135, 12).Value = 2250001620  # L135: was 9000000000
$ws.Cells.Item(135, 13).Value = -409106469  # M135: was -281260191
$ws.Cells.Item(135, 14).Value = -2250006690  # N135: was -9000005070
$ws.Cells.Item(138, 8).Value = 8976.869000000001  # H138: was 4098.6777
$ws.Cells.Item(138, 9).Value = 0  # I138: was 1141.0303
$ws.Cells.Item(138, 10).Value = 8976.869000000001  # J138: was 7852.615
$ws.Cells.Item(138, 11).Value = 0  # K138: was 3423.0909
$ws.Cells.Item(138, 12).Value = 26930.607  # L138: was 23557.845
$ws.Cells.Item(138, 13).ClearContents()  # M138: was 1716.9091
$ws.Cells.Item(138, 14).Value = -37210.607  # N138: was -33837.845
$ws.Cells.Item(141, 8).Value = 3119.5454  # H141: was 3250.647
$ws.Cells.Item(141, 9).Value = 1513.125  # I141: was 1837.5834
$ws.Cells.Item(141, 10).Value = 7403.3335  # J141: was 6642
$ws.Cells.Item(141, 11).Value = 4539.375  # K141: was 5512.7502
$ws.Cells.Item(141, 12).Value = 22210.0005  # L141: was 19926
$ws.Cells.Item(141, 13).Value = 640.625  # M141: was -332.7502000000004
$ws.Cells.Item(141, 14).Value = -32570.0005  # N141: was -30286

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 22683.773  # H32: was 21029.55
$ws.Cells.Item(32, 9).Value = 26307.096  # I32: was 23880.256
$ws.Cells.Item(32, 11).Value = 26307.096  # K32: was 23880.256
$ws.Cells.Item(32, 13).Value = -26020.096  # M32: was -23593.256
$ws.Cells.Item(34, 8).Value = 62016.8  # H34: was 0
$ws.Cells.Item(34, 9).Value = 50000  # I34: was 0
$ws.Cells.Item(34, 10).Value = 70028  # J34: was 0
$ws.Cells.Item(34, 11).Value = 50000  # K34: was 0
$ws.Cells.Item(34, 12).Value = 70028  # L34: was 0
$ws.Cells.Item(34, 13).Value = -49729  # M34: was None
$ws.Cells.Item(34, 14).Value = -70570  # N34: was None
$ws.Cells.Item(43, 8).Value = 22661  # H43: was 0
$ws.Cells.Item(43, 9).Value = 10342  # I43: was 0
$ws.Cells.Item(43, 10).Value = 34980  # J43: was 0
$ws.Cells.Item(43, 11).Value = 10342  # K43: was 0
$ws.Cells.Item(43, 12).Value = 34980  # L43: was 0
$ws.Cells.Item(43, 13).Value = -10029  # M43: was None
$ws.Cells.Item(43, 14).Value = -35606  # N43: was None
$ws.Cells.Item(61, 8).Value = 5827.8125  # H61: was 5288.7
$ws.Cells.Item(61, 9).Value = 4594.35  # I61: was 4120.9614
$ws.Cells.Item(61, 10).Value = 7883.5835  # J61: was 7457.357
$ws.Cells.Item(61, 11).Value = 4594.35  # K61: was 4120.9614
$ws.Cells.Item(61, 12).Value = 7883.5835  # L61: was 7457.357
$ws.Cells.Item(61, 13).Value = -4382.35  # M61: was -3908.9614
$ws.Cells.Item(61, 14).Value = -8307.583500000001  # N61: was -7881.357
$ws.Cells.Item(136, 8).Value = 5827.8125  # H136: was 5288.7
$ws.Cells.Item(136, 9).Value = 4594.35  # I136: was 4120.9614
$ws.Cells.Item(136, 10).Value = 7883.5835  # J136: was 7457.357
$ws.Cells.Item(136, 11).Value = 13783.05  # K136: was 12362.8842
$ws.Cells.Item(136, 12).Value = 23650.7505  # L136: was 22372.071
$ws.Cells.Item(136, 13).Value = -11233.05  # M136: was -9812.8842
$ws.Cells.Item(136, 14).Value = -28750.7505  # N136: was -27472.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(36, 8).Value = 850  # H36: was 2925
$ws.Cells.Item(36, 9).Value = 850  # I36: was 900
$ws.Cells.Item(36, 10).Value = 0  # J36: was 9000
$ws.Cells.Item(36, 11).Value = 850  # K36: was 900
$ws.Cells.Item(36, 12).Value = 0  # L36: was 9000
$ws.Cells.Item(36, 13).Value = -316  # M36: was -366
$ws.Cells.Item(36, 14).ClearContents()  # N36: was -10068
$ws.Cells.Item(130, 8).Value = 40000  # H130: was 60000
$ws.Cells.Item(130, 10).Value = 40000  # J130: was 60000
$ws.Cells.Item(130, 12).Value = 40000  # L130: was 60000
$ws.Cells.Item(130, 14).Value = -50040  # N130: was -70040
$ws.Cells.Item(134, 8).Value = 1910.7333  # H134: was 1192.9259
$ws.Cells.Item(134, 9).Value = 1422.3636  # I134: was 943.2174
$ws.Cells.Item(134, 10).Value = 3253.75  # J134: was 2628.75
$ws.Cells.Item(134, 11).Value = 4267.0908  # K134: was 2829.6522
$ws.Cells.Item(134, 12).Value = 9761.25  # L134: was 7886.25
$ws.Cells.Item(134, 13).Value = -1732.0908  # M134: was -294.6522
$ws.Cells.Item(134, 14).Value = -14831.25  # N134: was -12956.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 276.30768  # H22: was 268.72726
$ws.Cells.Item(22, 9).Value = 250.16667  # I22: was 243.33333
$ws.Cells.Item(22, 10).Value = 298.7143  # J22: was 299.2
$ws.Cells.Item(22, 11).Value = 250.16667  # K22: was 243.33333
$ws.Cells.Item(22, 12).Value = 298.7143  # L22: was 299.2
$ws.Cells.Item(22, 13).Value = 99.83332999999999  # M22: was 106.66667
$ws.Cells.Item(22, 14).Value = -998.7143  # N22: was -999.2
$ws.Cells.Item(31, 8).Value = 2665.257  # H31: was 2866.8125
$ws.Cells.Item(31, 9).Value = 1778.5555  # I31: was 1957.8667
$ws.Cells.Item(31, 10).Value = 3604.1177  # J31: was 3668.8235
$ws.Cells.Item(31, 11).Value = 1778.5555  # K31: was 1957.8667
$ws.Cells.Item(31, 12).Value = 3604.1177  # L31: was 3668.8235
$ws.Cells.Item(31, 13).Value = -1483.5555  # M31: was -1662.8667
$ws.Cells.Item(31, 14).Value = -4194.1177  # N31: was -4258.8235
$ws.Cells.Item(34, 8).Value = 2665.257  # H34: was 2866.8125
$ws.Cells.Item(34, 9).Value = 1778.5555  # I34: was 1957.8667
$ws.Cells.Item(34, 10).Value = 3604.1177  # J34: was 3668.8235
$ws.Cells.Item(34, 11).Value = 1778.5555  # K34: was 1957.8667
$ws.Cells.Item(34, 12).Value = 3604.1177  # L34: was 3668.8235
$ws.Cells.Item(34, 13).Value = -1576.5555  # M34: was -1755.8667
$ws.Cells.Item(34, 14).Value = -4008.1177  # N34: was -4072.8235
$ws.Cells.Item(94, 8).Value = 1039  # H94: was 1043.8235
$ws.Cells.Item(94, 10).Value = 1044.0588  # J94: was 1050.2
$ws.Cells.Item(94, 12).Value = 1044.0588  # L94: was 1050.2
$ws.Cells.Item(94, 14).Value = -1946.0588  # N94: was -1952.2
$ws.Cells.Item(105, 8).Value = 1542.3684  # H105: was 1636.7059
$ws.Cells.Item(105, 9).Value = 1180.3334  # I105: was 1274.6428
$ws.Cells.Item(105, 10).Value = 2900  # J105: was 3326.3333
$ws.Cells.Item(105, 11).Value = 1180.3334  # K105: was 1274.6428
$ws.Cells.Item(105, 12).Value = 2900  # L105: was 3326.3333
$ws.Cells.Item(105, 13).Value = 566.6666  # M105: was 472.3571999999999
$ws.Cells.Item(105, 14).Value = -6394  # N105: was -6820.3333
$ws.Cells.Item(122, 8).Value = 3934.0435  # H122: was 8474.682000000001
$ws.Cells.Item(122, 9).Value = 4720.1875  # I122: was 4531.1177
$ws.Cells.Item(122, 10).Value = 2137.1428  # J122: was 21882.8
$ws.Cells.Item(122, 11).Value = 14160.5625  # K122: was 13593.3531
$ws.Cells.Item(122, 12).Value = 6411.428400000001  # L122: was 65648.39999999999
$ws.Cells.Item(122, 13).Value = -11710.5625  # M122: was -11143.3531
$ws.Cells.Item(122, 14).Value = -11311.4284  # N122: was -70548.39999999999
$ws.Cells.Item(134, 8).Value = 3787.0667  # H134: was 3508.4443
$ws.Cells.Item(134, 9).Value = 3516.8333  # I134: was 2583.8823
$ws.Cells.Item(134, 10).Value = 4868  # J134: was 4335.684
$ws.Cells.Item(134, 11).Value = 10550.4999  # K134: was 7751.646900000001
$ws.Cells.Item(134, 12).Value = 14604  # L134: was 13007.052
$ws.Cells.Item(134, 13).Value = -8015.499899999999  # M134: was -5216.646900000001
$ws.Cells.Item(134, 14).Value = -19674  # N134: was -18077.052

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(14, 8).Value = 44300.715  # H14: was 37256.32
$ws.Cells.Item(14, 9).Value = 44300.715  # I14: was 37256.32
$ws.Cells.Item(14, 11).Value = 132902.145  # K14: was 111768.96
$ws.Cells.Item(14, 13).Value = -132729.145  # M14: was -111595.96
$ws.Cells.Item(113, 8).Value = 672.9231  # H113: was 657.6491
$ws.Cells.Item(113, 9).Value = 702.7778  # I113: was 697.2432
$ws.Cells.Item(113, 10).Value = 605.75  # J113: was 584.4
$ws.Cells.Item(113, 11).Value = 2108.3334  # K113: was 2091.7296
$ws.Cells.Item(113, 12).Value = 1817.25  # L113: was 1753.2
$ws.Cells.Item(113, 13).Value = 61.66660000000002  # M113: was 78.27039999999988
$ws.Cells.Item(113, 14).Value = -6157.25  # N113: was -6093.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 764.04  # H93: was 883.3182
$ws.Cells.Item(93, 9).Value = 492.83334  # I93: was 515.7646999999999
$ws.Cells.Item(93, 10).Value = 1461.4286  # J93: was 2133
$ws.Cells.Item(93, 11).Value = 492.83334  # K93: was 515.7646999999999
$ws.Cells.Item(93, 12).Value = 1461.4286  # L93: was 2133
$ws.Cells.Item(93, 13).Value = 755.16666  # M93: was 732.2353000000001
$ws.Cells.Item(93, 14).Value = -3957.4286  # N93: was -4629
$ws.Cells.Item(97, 8).Value = 28221.75  # H97: was 33959
$ws.Cells.Item(97, 10).Value = 28221.75  # J97: was 33959
$ws.Cells.Item(97, 12).Value = 28221.75  # L97: was 33959
$ws.Cells.Item(97, 14).Value = -30203.75  # N97: was -35941
$ws.Cells.Item(116, 8).Value = 40000  # H116: was 0
$ws.Cells.Item(116, 10).Value = 40000  # J116: was 0
$ws.Cells.Item(116, 12).Value = 40000  # L116: was 0
$ws.Cells.Item(116, 14).Value = -49178  # N116: was None
$ws.Cells.Item(132, 8).Value = 5175.7036  # H132: was 5490.0435
$ws.Cells.Item(132, 9).Value = 4930.0557  # I132: was 5376.2856
$ws.Cells.Item(132, 11).Value = 14790.1671  # K132: was 16128.8568
$ws.Cells.Item(132, 13).Value = -12260.1671  # M132: was -13598.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2256.7  # H132: was 2087.9092
$ws.Cells.Item(132, 9).Value = 1082.6875  # I132: was 974.8946999999999
$ws.Cells.Item(132, 11).Value = 3248.0625  # K132: was 2924.6841
$ws.Cells.Item(132, 13).Value = -718.0625  # M132: was -394.6840999999999
